# Apply the "biochar_land" scenario split + plant-cost updates described in
# the commit:
#   "removed bio-oil production in biochar land scenario and updated plant costs"
#
# 1. Rename the existing sheet ("Sheet1" -> "pyrolysis-nofert").
# 2. Add the new shared string used at the bottom of both sheets
#    ("same as beef/dairy") before any other new string, so it lands at the
#    same sharedStrings index the target workbook uses.
# 3. Append a "same as beef/dairy" note in column B of the last data row of
#    the pyrolysis-nofert sheet.
# 4. Duplicate pyrolysis-nofert into a new "biochar_land" sheet (this keeps
#    every formula / number style identical to start with).
# 5. In biochar_land: drop the "GJ/kg feedstock" row (row 9) entirely (the
#    biochar-land scenario no longer tracks bio-oil yield from that input),
#    which also shifts every row below it up by one and re-points the
#    dependent formulas automatically.
# 6. In biochar_land: remove the now-unused "Yield (kg/GJ)" column (J) and
#    blank out the "Yield (GJ/kg feedstock)" column (I) -- the biochar-land
#    scenario no longer produces bio-oil, so those derived yield figures go
#    away, but column I keeps its currency number format for future use.
# 7. In biochar_land: replace the hard-coded bio-oil yield fractions in
#    A4:A8 with the literal "1/<moisture ratio>" formulas used for the
#    biochar-land plant-cost recompute.
# 8. In biochar_land: relabel the unit-cost header from "$1975/GJ" to
#    "$1975/kg" (cost is now tracked per kg, not per GJ of bio-oil).
# 9. Fix up tab/selection state: biochar_land ends with row 9 selected
#    (reflecting the deleted "GJ/kg feedstock" row), pyrolysis-nofert ends
#    active with E30 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- pyrolysis-nofert (sheet 1) -------------------------------------------
$ws1.Name = "pyrolysis-nofert"

# Add the note first so "same as beef/dairy" gets the lower shared-string
# index (it is referenced from both sheets).
$ws1.Range("B26").Value = "same as beef/dairy"

# --- biochar_land (sheet 2, duplicated from pyrolysis-nofert) -------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "biochar_land"

# Drop the "GJ/kg feedstock" row; everything below shifts up and formulas
# (A$10 -> A$9, etc.) re-point automatically.
$ws2.Rows("9:9").Delete()

# The "Yield (kg/GJ)" column is removed outright.
$ws2.Columns("J:J").Delete()

# The "Yield (GJ/kg feedstock)" column's formulas are cleared (biochar_land
# no longer produces bio-oil), but the column/cell formatting stays.
$ws2.Range("I12:I17").ClearContents()

# Bio-oil yield fractions become explicit "1/x" formulas for biochar_land.
$ws2.Range("A4").Formula = "=1/2.1815"
$ws2.Range("A5").Formula = "=1/2.1052"
$ws2.Range("A6").Formula = "=1/2.055"
$ws2.Range("A7").Formula = "=1/2.136"
$ws2.Range("A8").Formula = "=1/2.1276"

# Unit-cost header is now per kg instead of per GJ.
$ws2.Range("H12").Value = "Unit cost  (`$1975/kg)"

# --- final selection / active-sheet state ----------------------------------
$ws2.Rows("9:9").Select()
$ws1.Activate()
$ws1.Range("E30").Select()

Write-Output "biochar_land scenario created; plant costs updated"
